# Update the "timestamp" column (Z) for the existing data rows (2-29) to
# reflect the latest run of the pcSMOTE logging pass - the synthetic
# samples were successfully plotted, so the log timestamps were refreshed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @{
    2  = "2025-11-02T02:04:09.106519"
    3  = "2025-11-02T02:04:09.106519"
    4  = "2025-11-02T02:04:09.106519"
    5  = "2025-11-02T02:04:09.106519"
    6  = "2025-11-02T02:04:09.106519"
    7  = "2025-11-02T02:04:09.106519"
    8  = "2025-11-02T02:04:09.106519"
    9  = "2025-11-02T02:04:09.107526"
    10 = "2025-11-02T02:04:09.107526"
    11 = "2025-11-02T02:04:09.107526"
    12 = "2025-11-02T02:04:09.108108"
    13 = "2025-11-02T02:04:09.108108"
    14 = "2025-11-02T02:04:09.108108"
    15 = "2025-11-02T02:04:09.108108"
    16 = "2025-11-02T02:04:09.108108"
    17 = "2025-11-02T02:04:09.108672"
    18 = "2025-11-02T02:04:09.108672"
    19 = "2025-11-02T02:04:09.108672"
    20 = "2025-11-02T02:04:09.108672"
    21 = "2025-11-02T02:04:09.108672"
    22 = "2025-11-02T02:04:09.109202"
    23 = "2025-11-02T02:04:09.109202"
    24 = "2025-11-02T02:04:09.109202"
    25 = "2025-11-02T02:04:09.109202"
    26 = "2025-11-02T02:04:09.109202"
    27 = "2025-11-02T02:04:09.109202"
    28 = "2025-11-02T02:04:09.109202"
    29 = "2025-11-02T02:04:09.109202"
}

foreach ($row in $timestamps.Keys) {
    $ws.Range("Z$row").Value = $timestamps[$row]
}
